$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.528410666666667
$ws.Range("H2").Value = 22.585232
$ws.Range("I2").Value = 0.3846034394416388
$ws.Range("J2").Value = 0.3846034394416388
$ws.Range("M2").Value = 1.599392
$ws.Range("N2").Value = 4.798176
$ws.Range("O2").Value = 0.03952976301548796
$ws.Range("P2").Value = 0.03952976301548796
$ws.Range("Q2").Value = 12.04087979298133
$ws.Range("R2").Value = 108.367918136832
$ws.Range("S2").Value = 0.01520328281606956
$ws.Range("T2").Value = 0.01520328281606956
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.528410666666667
$ws.Range("H3").Value = 22.585232
$ws.Range("I3").Value = 0.3846034394416388
$ws.Range("J3").Value = 0.3846034394416388
$ws.Range("O3").Value = 0.4638329693976876
$ws.Range("P3").Value = 0.4638329693976876
$ws.Range("Q3").Value = 141.2848598751005
$ws.Range("R3").Value = 1271.563738875904
$ws.Range("S3").Value = 0.1783917553567791
$ws.Range("T3").Value = 0.1783917553567791
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.528410666666667
$ws.Range("H4").Value = 22.585232
$ws.Range("I4").Value = 0.3846034394416388
$ws.Range("J4").Value = 0.3846034394416388
$ws.Range("M4").Value = 20.09416733333333
$ws.Range("N4").Value = 60.28250199999999
$ws.Range("O4").Value = 0.4966372675868244
$ws.Range("P4").Value = 0.4966372675868245
$ws.Range("Q4").Value = 151.2771436900515
$ws.Range("R4").Value = 1361.494293210464
$ws.Range("S4").Value = 0.1910084012687902
$ws.Range("T4").Value = 0.1910084012687902
$ws.Range("I5").Value = 0.4334462803064951
$ws.Range("J5").Value = 0.433446280306495
$ws.Range("M5").Value = 1.599392
$ws.Range("N5").Value = 4.798176
$ws.Range("O5").Value = 0.03952976301548796
$ws.Range("P5").Value = 0.03952976301548796
$ws.Range("Q5").Value = 13.57001530059733
$ws.Range("R5").Value = 122.130137705376
$ws.Range("S5").Value = 0.01713402874046051
$ws.Range("T5").Value = 0.01713402874046051
$ws.Range("I6").Value = 0.4334462803064951
$ws.Range("J6").Value = 0.433446280306495
$ws.Range("O6").Value = 0.4638329693976876
$ws.Range("P6").Value = 0.4638329693976876
$ws.Range("S6").Value = 0.201046675268944
$ws.Range("T6").Value = 0.201046675268944
$ws.Range("I7").Value = 0.4334462803064951
$ws.Range("J7").Value = 0.433446280306495
$ws.Range("M7").Value = 20.09416733333333
$ws.Range("N7").Value = 60.28250199999999
$ws.Range("O7").Value = 0.4966372675868244
$ws.Range("P7").Value = 0.4966372675868245
$ws.Range("Q7").Value = 170.4886345349336
$ws.Range("R7").Value = 1534.397710814402
$ws.Range("S7").Value = 0.2152655762970905
$ws.Range("T7").Value = 0.2152655762970905
$ws.Range("G8").Value = 3.561581333333333
$ws.Range("H8").Value = 10.684744
$ws.Range("I8").Value = 0.1819502802518661
$ws.Range("J8").Value = 0.1819502802518661
$ws.Range("M8").Value = 1.599392
$ws.Range("N8").Value = 4.798176
$ws.Range("O8").Value = 0.03952976301548796
$ws.Range("P8").Value = 0.03952976301548796
$ws.Range("Q8").Value = 5.696364691882667
$ws.Range("R8").Value = 51.267282226944
$ws.Range("S8").Value = 0.007192451458957885
$ws.Range("T8").Value = 0.007192451458957885
$ws.Range("G9").Value = 3.561581333333333
$ws.Range("H9").Value = 10.684744
$ws.Range("I9").Value = 0.1819502802518661
$ws.Range("J9").Value = 0.1819502802518661
$ws.Range("O9").Value = 0.4638329693976876
$ws.Range("P9").Value = 0.4638329693976876
$ws.Range("Q9").Value = 66.83980748310756
$ws.Range("R9").Value = 601.5582673479681
$ws.Range("S9").Value = 0.08439453877196448
$ws.Range("T9").Value = 0.08439453877196448
$ws.Range("G10").Value = 3.561581333333333
$ws.Range("H10").Value = 10.684744
$ws.Range("I10").Value = 0.1819502802518661
$ws.Range("J10").Value = 0.1819502802518661
$ws.Range("M10").Value = 20.09416733333333
$ws.Range("N10").Value = 60.28250199999999
$ws.Range("O10").Value = 0.4966372675868244
$ws.Range("P10").Value = 0.4966372675868245
$ws.Range("Q10").Value = 71.56701128327643
$ws.Range("R10").Value = 644.103101549488
$ws.Range("S10").Value = 0.0903632900209437
$ws.Range("T10").Value = 0.09036329002094372
